# Comentado el agregar imagen en el formulario de humedal para que no aparezca
# Adds "Mi Rama" tracking section (rows 14-25) with required-field warning
# messages for Cuenca, Complejo, Persona, Presiones, Fauna and Flora.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Clear out the old rows 14 and 15 completely so we can rebuild the
#    block from scratch (they get replaced by the new rows 14-25).
# ---------------------------------------------------------------------
$ws.Range("A14:I15").Clear()

# ---------------------------------------------------------------------
# 2. Row 14 - new header-like row ("Pruebas" / "Mi Rama" / "Rama Principal (BD Nueva)")
#    Reuses the formatting of row 1's header cells.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)

$ws.Range("D1").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)

$ws.Range("F1").Copy()
$ws.Range("E14").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("C14").Value = "Pruebas"
$ws.Range("D14").Value = "Mi Rama"
$ws.Range("E14").Value = "Rama Principal (BD Nueva)"

$ws.Rows.Item(14).RowHeight = 60.75

# ---------------------------------------------------------------------
# 3. Rows 15-20 - warning-message rows (with label text in column C)
# ---------------------------------------------------------------------
$labelSourceRows = @{
    15 = "C2"
    16 = "C3"
    17 = "C3"
    18 = "C5"
    19 = "C3"
    20 = "C3"
}
$labelTexts = @{
    15 = "Mensaje de adbertencia campo obligatorio incompleto (Cuenca)"
    16 = "Mensaje de adbertencia campo obligatorio incompleto (Complejo)"
    17 = "Mensaje de adbertencia campo obligatorio incompleto (Persona)"
    18 = "Mensaje de adbertencia campo obligatorio incompleto (Presiones)"
    19 = "Mensaje de adbertencia campo obligatorio incompleto (Fauna)"
    20 = "Mensaje de adbertencia campo obligatorio incompleto (Flora)"
}

foreach ($r in 15..20) {
    $src = $labelSourceRows[$r]

    $ws.Range($src).Copy()
    $ws.Range("C$r").PasteSpecial($xlPasteFormats)

    $ws.Range("D2").Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteFormats)

    $ws.Range("F2").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("C$r").WrapText = $true
    $ws.Range("C$r").Value = $labelTexts[$r]
    $ws.Range("D$r").Value = [char]0x2714
    $ws.Range("E$r").Value = "O"

    $ws.Rows.Item($r).RowHeight = 30.75
}

# ---------------------------------------------------------------------
# 4. Rows 21-25 - same pattern but column C stays empty
# ---------------------------------------------------------------------
$emptySourceRows = @{
    21 = "C3"
    22 = "C3"
    23 = "C5"
    24 = "C3"
    25 = "C3"
}

foreach ($r in 21..25) {
    $src = $emptySourceRows[$r]

    $ws.Range($src).Copy()
    $ws.Range("C$r").PasteSpecial($xlPasteFormats)

    $ws.Range("D2").Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteFormats)

    $ws.Range("F2").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("C$r").WrapText = $true
    $ws.Range("D$r").Value = [char]0x2714
    $ws.Range("E$r").Value = "O"

    $ws.Rows.Item($r).RowHeight = 29.25
}

# ---------------------------------------------------------------------
# 5. Make sure column F has nothing in the new rows (it is unused there)
# ---------------------------------------------------------------------
$ws.Range("F14:F25").Clear()

# ---------------------------------------------------------------------
# 6. Update the view: selection moves to G21, no explicit scroll anchor
# ---------------------------------------------------------------------
$ws.Range("G21").Select()
